# Track positive expenses and negative credits
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# --- Rename the tax_category header to a friendlier label, with wrap text ---
$ws.Range("E1").Value = "Tax category"
$ws.Range("E1").WrapText = $true

# --- Format the "Another extra" column (G) with an example date, like the
#     other "Eg ..." placeholder cells in row 2 ---
$ws.Range("G1:G2").NumberFormat = "m/d/yy;@"
$ws.Cells.Item(2, 7).Value = 43466

# --- Give the Date column (A) a proper date display instead of a raw
#     serial number ---
$ws.Range("A1:A5").NumberFormat = "m/d/yy"

# --- Add two new transactions: a positive expense and a negative credit ---
$ws.Cells.Item(6, 4).Value = 101
$ws.Cells.Item(6, 5).Value = "Other"
$ws.Cells.Item(7, 4).Value = -100
$ws.Cells.Item(7, 5).Value = "Other"

# --- Widen the amount column slightly and select the newest entry ---
$ws.Columns.Item(4).ColumnWidth = 7.5
$ws.Range("D7").Select()

# --- Zoom in a bit on the Transactions sheet ---
$win = $excel.ActiveWindow
$win.Zoom = 130

# --- Add a new, empty worksheet after Transactions ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$newSheet.Name = "Extra sheet"
